$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2022" data column (S), mirroring the formatting used in
# the existing "2021" column (R) for each of the three populated rows.
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("S4").Value = 2022

$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("S5").Value = 30

$ws.Range("R6").Copy()
$ws.Range("S6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("S6").Value = 11928.6

$excel.CutCopyMode = $false

# Move the active selection as recorded after the edit.
$ws.Range("T3").Select()
